$wb = $excel.ActiveWorkbook

# Rename the *img sheets to img* (reverse prefix/suffix order)
$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# Move the active tab from "holiday" to the renamed "imge" sheet (was "eimg")
$wb.Worksheets.Item("imge").Activate()
